$d = $word.ActiveDocument

# --- Update existing paragraphs 1-6 in place (text only; keep paragraph marks) ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.MoveEnd(1, -1) | Out-Null
$r1.Text = "⚡️🚀המאמר היומי של מייק 09.07.24: ⚡️🚀"

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.MoveEnd(1, -1) | Out-Null
$r2.Text = "Learning to (Learn at Test Time): RNNs with Expressive Hidden States"

$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$r3.MoveEnd(1, -1) | Out-Null
$r3.Text = "המאמר הזה המצהיר שהוא לומד ב״זמן טסט״ משך את עיניי היום. המאמר מציע ארכיטקטורה חדשה ומעניינת לעיבוד דאטה סדרתי. בעיקרון הרשת די דומה ל-RNN מבחינת המהות אבל יש כמה הבדלים מהותיים."

$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.MoveEnd(1, -1) | Out-Null
$r4.Text = "ֿאז מה יש לנו בארכיטקטורה הזו? בדומה ל-RNN אנו מחשבים את הייצוג עבור יחידת דאטה בזמן t (נגיד טוקן t) אבל כאן עושים זאת בשיטה שונה. לפי המאמר במקום לחשב את הייצוג עצמו אנו מחשבים את וקטור המשקלים שיאפשר לנו לחשב את ייצוגו של יחידת דאטה t. כלומר אנו מעדכנים את משקלות מודל בתנועה בהתאם לדאטה כלומר הרשת מתאפטמת ומתאימה את עצמה לדאטה שעליה היא מופעלת. זה נעשה באמצעות הזזה של המקשלים בכיוון הנגדי של הגרדיאנט של פונקציית לוס l."

$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5.MoveEnd(1, -1) | Out-Null
$r5.Text = "מה זה בעצם פונקציית l ואיך מאמנים אותה? נניח שהייצוג של איבר דאטה t מחושב על ידי פונקציית f. במקרה הזה פונקציית l יכולה להיות (למשל) נורמה של הפרש ריבוע של ייצוג דאטה z (המחושב עם f) מהדאטה עצמו. כלומר אנו מאמנים את וקטור הייצוג להיות מסוגל לשחזר (כלומר לזכור) את הדאטה עצמו x_t. כמובן שאין בזה הרבה משמעות אבל אם נאמן רשת עם קלט מורעש ונשווה את ייצוג עם הדאטה האמיתי נקבל סוג של רשת denoising שהרשת לומדת להפיק ייצוג המאפשר לזכור את הפיצ'רים המהותיים של דאטה הנחוצים לשחזור."

$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
$r6.MoveEnd(1, -1) | Out-Null
$r6.Text = "דרך אחרת המוצעת במאמר לאמן את רשת לשחזר הטלה למימד נמוך של דאטה להטלה אחרת כאשר שתי ההטלות נלמדות גם כן. הייצוג של דאטה במקרה הזה מחושב עם הטלה נלמדת שלישית (עם פונקציית f). כלומר המטרה כאן ללמוד את ייצוג של דאטה כאשר המשקלים מחושבים עם GD מהמשקלים הקודמים."

# --- Insert two new paragraphs after paragraph 6 (before the link paragraph content) ---
$p6 = $d.Paragraphs.Item(6)
$newPara1 = $p6.Range.InsertParagraphAfter()
$d.Paragraphs.Item(7).Range.Text = "הארכיטקטורה קיבלה שם ttt וניתן לשלב אותם על שכבות אחרות (כמו טרנספורמרים או SSM). רעיון מגניב שבינתיים לא הפנמתי אותו עד הסוף…"

$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(8).Range.Text = "https://arxiv.org/pdf/2407.04620"

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
